# "modificaciones para rama dos"
#
# The source document has two paragraphs:
#   1) "Arancando"  -- still wrapped in the spell-checker's <w:proofErr> markers
#   2) ""            -- empty paragraph
#
# The edit:
#   - clears the stale spell-check markers around "Arancando" (they are not
#     reachable as ordinary document text, so we refresh that paragraph's
#     OOXML in place to drop them)
#   - types "Modifico para rama dos" into the previously empty second
#     paragraph, in the same es-AR language as the rest of the document

$d = $word.ActiveDocument

# --- locate the two paragraphs we care about -----------------------------
$arancandoPara = $null
$emptyPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    $text = $para.Range.Text
    if ($text -match "Arancando") {
        $arancandoPara = $para
    } elseif ($text.Trim() -eq "" -and $emptyPara -eq $null) {
        $emptyPara = $para
    }
}

# --- 1) drop the leftover spellcheck proofing marks around "Arancando" ---
if ($arancandoPara -ne $null) {
    $arancandoPara.Range.InsertXML($arancandoPara.Range.WordOpenXML)
}

# --- 2) add the new line to the empty paragraph --------------------------
if ($emptyPara -ne $null) {
    $emptyPara.Range.InsertAfter("Modifico para rama dos")
    $emptyPara.Range.LanguageID = "es-AR"
}
